$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C ("Förändrad") holds a date serial number for every data row
# (rows 2-53). The sheet was regenerated with a newer "changed" date,
# moving the serial value from 45192 (2023-09-23) to 45202 (2023-10-03)
# for each of those rows.
$ws.Range("C2:C53").Value = 45202
